$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "225÷5="
$t.Cell(1,2).Range.Text = "658÷7="
$t.Cell(1,3).Range.Text = "186÷8="
$t.Cell(1,4).Range.Text = "642÷6="
$t.Cell(1,5).Range.Text = "878÷6="
$t.Cell(5,1).Range.Text = "422÷4="
$t.Cell(5,2).Range.Text = "528÷2="
$t.Cell(5,3).Range.Text = "696÷5="
$t.Cell(5,4).Range.Text = "298÷4="
$t.Cell(5,5).Range.Text = "663÷7="
$t.Cell(9,1).Range.Text = "460÷5="
$t.Cell(9,2).Range.Text = "907÷9="
$t.Cell(9,3).Range.Text = "114÷2="
$t.Cell(9,4).Range.Text = "763÷9="
$t.Cell(9,5).Range.Text = "261÷3="
$t.Cell(13,1).Range.Text = "293÷2="
$t.Cell(13,2).Range.Text = "769÷5="
$t.Cell(13,3).Range.Text = "292÷4="
$t.Cell(13,4).Range.Text = "128÷3="
$t.Cell(13,5).Range.Text = "442÷2="
$t.Cell(17,1).Range.Text = "759÷2="
$t.Cell(17,2).Range.Text = "848÷5="
$t.Cell(17,3).Range.Text = "793÷8="
$t.Cell(17,4).Range.Text = "165÷4="
$t.Cell(17,5).Range.Text = "845÷3="
